$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09124040616993714
$ws.Range("H2").Value = -30.57381586398508
$ws.Range("I2").Value = 1.769845217231135
$ws.Range("G3").Value = 0.147159293691548
$ws.Range("H3").Value = 65.27844520469607
$ws.Range("G4").Value = -0.7945873311133941
$ws.Range("H4").Value = -29.12908218929165
$ws.Range("G5").Value = -0.6515421682246798
$ws.Range("H5").Value = -6.706416064839506
$ws.Range("G6").Value = 0.1996875845601599
$ws.Range("H6").Value = -18.86173541420811
$ws.Range("G7").Value = 0.3882072456746861
$ws.Range("H7").Value = 136.9690159701107
$ws.Range("G8").Value = 0.1594138041217252
$ws.Range("H8").Value = -3.566125750523551
$ws.Range("G9").Value = 0.2077918103965916
$ws.Range("H9").Value = 6.509357851903191
$ws.Range("G10").Value = -0.1238758092119051
$ws.Range("H10").Value = -116.7694250701831
$ws.Range("G11").Value = -0.1069914105843795
$ws.Range("H11").Value = 9.915869935000034
$ws.Range("G12").Value = 0.2089564476610566
$ws.Range("H12").Value = 31.39494818062817
$ws.Range("G13").Value = 0.176582253717831
$ws.Range("H13").Value = -14.13947608827334
$ws.Range("G14").Value = 0.1783164488523147
$ws.Range("H14").Value = -5.838542454871688
$ws.Range("G15").Value = 0.2300177221947939
$ws.Range("H15").Value = -7.958088949487803
$ws.Range("G16").Value = -0.006705822848939855
$ws.Range("H16").Value = -118.3820654183642
$ws.Range("G17").Value = -0.0001456530265392189
$ws.Range("H17").Value = -100.4106304020959
$ws.Range("G18").Value = 0.04465794095562786
$ws.Range("H18").Value = -74.23296248552464
$ws.Range("G19").Value = 0.06007550535521902
$ws.Range("H19").Value = -52.22735666143728
$ws.Range("G20").Value = 0.07899718956387268
$ws.Range("H20").Value = -31.09828959952467
$ws.Range("G21").Value = 0.09526700080846671
$ws.Range("H21").Value = -5.111439386570639
$ws.Range("G22").Value = 0.06423643509447831
$ws.Range("H22").Value = -31.80637703705857
$ws.Range("G23").Value = 0.08041847627501268
$ws.Range("H23").Value = -25.87474857737158
$ws.Range("G24").Value = -0.1971171139517948
$ws.Range("H24").Value = -58.16902401818408
$ws.Range("G25").Value = -0.1398674342544872
$ws.Range("H25").Value = 37.12525747910595
$ws.Range("G26").Value = 0.2250979697836644
$ws.Range("H26").Value = 41.58724124878993
$ws.Range("G27").Value = 0.1719898235344214
$ws.Range("H27").Value = -14.20513458211327
$ws.Range("G28").Value = 0.0305621861390063
$ws.Range("H28").Value = 480.5787388926575
$ws.Range("G29").Value = 0.01311046931265
$ws.Range("H29").Value = -14.74247266651109
